# Fix excel. Mobile login
# - Add a new "NIVEL" column (C) to the PREGUNTAS sheet with level values
#   for each question, formatted like the rest of the table.
# - Make PREGUNTAS the active/selected sheet (instead of RESPUESTAS).

$wb = $excel.ActiveWorkbook

$preguntas = $wb.Worksheets.Item("PREGUNTAS")

# --- New header cell C1 = "NIVEL", matching the look of the existing
#     header cells (A1 "PREGUNTA" / B1 "COMENTARIO").
$preguntas.Range("A1").Copy()
$preguntas.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$preguntas.Range("C1").Value = "NIVEL"

# --- New data values for the NIVEL column, centered like the rest of the
#     data rows.
$preguntas.Range("C2").Value = 1
$preguntas.Range("C3").Value = 2
$preguntas.Range("C2:C3").HorizontalAlignment = -4108

# --- Make PREGUNTAS the active sheet/tab and set its selection.
$preguntas.Select() | Out-Null
$preguntas.Range("C4").Select() | Out-Null
